{"js": "// Update each three-digit x one-digit multiplication answer in the table\n// to its new value, matched by searching for the old expression text.\nconst body = context.document.body;\n\nconst replacements = [\n  [\"461\u00d72=922\", \"624\u00d79=5616\"],\n  [\"643\u00d76=3858\", \"523\u00d77=3661\"],\n  [\"789\u00d73=2367\", \"550\u00d77=3850\"],\n  [\"115\u00d79=1035\", \"466\u00d77=3262\"],\n  [\"433\u00d73=1299\", \"452\u00d72=904\"],\n  [\"308\u00d78=2464\", \"238\u00d72=476\"],\n  [\"494\u00d73=1482\", \"732\u00d78=5856\"],\n  [\"470\u00d73=1410\", \"486\u00d75=2430\"],\n  [\"697\u00d75=3485\", \"260\u00d75=1300\"],\n  [\"177\u00d74=708\", \"278\u00d79=2502\"],\n  [\"583\u00d76=3498\", \"744\u00d74=2976\"],\n  [\"598\u00d78=4784\", \"822\u00d76=4932\"],\n  [\"226\u00d79=2034\", \"273\u00d77=1911\"],\n  [\"359\u00d75=1795\", \"539\u00d78=4312\"],\n  [\"651\u00d76=3906\", \"796\u00d78=6368\"],\n  [\"657\u00d76=3942\", \"843\u00d72=1686\"],\n  [\"153\u00d73=459\", \"357\u00d75=1785\"],\n  [\"957\u00d72=1914\", \"852\u00d76=5112\"],\n  [\"614\u00d79=5526\", \"464\u00d75=2320\"],\n  [\"714\u00d77=4998\", \"727\u00d74=2908\"],\n  [\"724\u00d75=3620\", \"698\u00d76=4188\"],\n  [\"951\u00d76=5706\", \"526\u00d77=3682\"],\n  [\"589\u00d73=1767\", \"725\u00d75=3625\"],\n  [\"116\u00d79=1044\", \"993\u00d76=5958\"],\n  [\"260\u00d76=1560\", \"748\u00d72=1496\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit multiplication answer with its updated value.\n# Uses Word Find/Replace (Content.Find) against the whole document body.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"461\u00d72=922\"\n$find.Replacement.Text = \"624\u00d79=5616\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"643\u00d76=3858\"\n$find.Replacement.Text = \"523\u00d77=3661\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"789\u00d73=2367\"\n$find.Replacement.Text = \"550\u00d77=3850\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"115\u00d79=1035\"\n$find.Replacement.Text = \"466\u00d77=3262\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"433\u00d73=1299\"\n$find.Replacement.Text = \"452\u00d72=904\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"308\u00d78=2464\"\n$find.Replacement.Text = \"238\u00d72=476\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"494\u00d73=1482\"\n$find.Replacement.Text = \"732\u00d78=5856\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"470\u00d73=1410\"\n$find.Replacement.Text = \"486\u00d75=2430\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"697\u00d75=3485\"\n$find.Replacement.Text = \"260\u00d75=1300\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"177\u00d74=708\"\n$find.Replacement.Text = \"278\u00d79=2502\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"583\u00d76=3498\"\n$find.Replacement.Text = \"744\u00d74=2976\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"598\u00d78=4784\"\n$find.Replacement.Text = \"822\u00d76=4932\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"226\u00d79=2034\"\n$find.Replacement.Text = \"273\u00d77=1911\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"359\u00d75=1795\"\n$find.Replacement.Text = \"539\u00d78=4312\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"651\u00d76=3906\"\n$find.Replacement.Text = \"796\u00d78=6368\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"657\u00d76=3942\"\n$find.Replacement.Text = \"843\u00d72=1686\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"153\u00d73=459\"\n$find.Replacement.Text = \"357\u00d75=1785\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"957\u00d72=1914\"\n$find.Replacement.Text = \"852\u00d76=5112\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"614\u00d79=5526\"\n$find.Replacement.Text = \"464\u00d75=2320\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"714\u00d77=4998\"\n$find.Replacement.Text = \"727\u00d74=2908\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"724\u00d75=3620\"\n$find.Replacement.Text = \"698\u00d76=4188\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"951\u00d76=5706\"\n$find.Replacement.Text = \"526\u00d77=3682\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"589\u00d73=1767\"\n$find.Replacement.Text = \"725\u00d75=3625\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"116\u00d79=1044\"\n$find.Replacement.Text = \"993\u00d76=5958\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"260\u00d76=1560\"\n$find.Replacement.Text = \"748\u00d72=1496\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n"}
